# Add Vuelta 2023, bar_plot and final confusion matrix
# - Correct a handful of stat columns on the first Vuelta stage row (40)
# - Correct nbr-3-cat / nbr-2-cat / elevation-gain on row 52
# - Populate the new "Breakaway win" (L) column for the Vuelta stage rows (40-58)
# - Move the active selection/scroll position to reflect the newly added data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40 corrections ---
$ws.Range("D40").Value = 182
$ws.Range("F40").Value = 3
$ws.Range("I40").Value = 2
$ws.Range("K40").Value = 2754

# --- Row 52 corrections ---
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 2
$ws.Range("K52").Value = 2494

# --- New "Breakaway win" (column L) values for Vuelta rows 40-58 ---
$breakawayWin = @{
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 0
    50 = 0
    51 = 1
    52 = 1
    53 = 0
    54 = 0
    55 = 1
    56 = 0
    57 = 1
    58 = 0
}

foreach ($row in $breakawayWin.Keys) {
    $ws.Cells.Item($row, 12).Value = $breakawayWin[$row]
}

# --- Move selection / scroll position to the bottom of the refreshed table ---
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("L59").Select()
